# Update status ("Estado") values for several requirements on Sheet1,
# clear out the now-obsolete standalone note in A15, resize columns to fit
# the new content, and leave the selection on D4 (last touched cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Status changes in column D (uses existing "Estados" list: Pendiente/En Progreso/Listo)
$ws.Range("D2").Value = "Listo"
$ws.Range("D3").Value = "Listo"
$ws.Range("D4").Value = "Listo"
$ws.Range("D5").Value = "Pendiente"
$ws.Range("D8").Value = "Listo"

# Remove the stray leftover note in A15 (row disappears once it's empty)
$ws.Range("A15").ClearContents()

# Re-fit the columns now that cell contents/lengths changed
$ws.Columns.Item(1).ColumnWidth = 2.0221354166666665
$ws.Columns.Item(2).ColumnWidth = 71.73697916666667
$ws.Columns.Item(4).ColumnWidth = 10.022135416666666

# Leave the active selection on D4
[void]$ws.Range("D4").Select()
